# Horarios actualizados Linea 141 - 485
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "LP1912": refresh scrape timestamp / row count, and update the data
# rows (6-20 changed, 21-22 newly appended).
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 05:15:38"
$ws1.Range("A3").Value = "Total filas: 17"

$data1 = @(
    @("05:15:38", "05:27", "23_HERNANDEZ", 12, "LP1912"),
    @("05:15:38", "05:31", "81_EL PELIGRO", 16, "LP1912"),
    @("05:15:38", "05:44", "14_ABASTO", 29, "LP1912"),
    @("05:15:38", "05:52", "17_ROMERO", 37, "LP1912"),
    @("05:15:38", "06:01", "16_SANTA ANA", 46, "LP1912"),
    @("05:15:38", "06:03", "10_OLMOS", 48, "LP1912"),
    @("05:15:38", "06:10", "215A_EL PATO", 55, "LP1912"),
    @("05:15:38", "06:24", "11_ETCHEVERRY", 69, "LP1912"),
    @("05:15:38", "06:27", "23_HERNANDEZ", 72, "LP1912"),
    @("05:15:38", "06:31", "17X38_ROMERO", 76, "LP1912"),
    @("05:15:38", "06:31", "16_SANTA ANA", 76, "LP1912"),
    @("05:15:38", "06:39", "225_C ROCA-H SUR", 84, "LP1912"),
    @("05:15:38", "06:50", "215A_EL PATO", 95, "LP1912"),
    @("05:15:38", "06:54", "14_ABASTO", 99, "LP1912"),
    @("05:15:38", "07:04", "225_GOMEZ", 109, "LP1912"),
    @("05:15:38", "07:06", "215C_EL PATO", 111, "LP1912"),
    @("05:15:38", "07:13", "14X44_ABASTO", 118, "LP1912")
)

$r = 6
foreach ($row in $data1) {
    $ws1.Cells.Item($r, 1).Value = $row[0]
    $ws1.Cells.Item($r, 2).Value = $row[1]
    $ws1.Cells.Item($r, 3).Value = $row[2]
    $ws1.Cells.Item($r, 4).Value = $row[3]
    $ws1.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet "LP1912-215": refresh scrape timestamp / row count, update existing
# rows and append the newly scraped row.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 05:15:38"
$ws2.Range("A3").Value = "Total filas: 3"

$data2 = @(
    @("05:15:38", "06:10", "215A_EL PATO", 55, "LP1912"),
    @("05:15:38", "06:50", "215A_EL PATO", 95, "LP1912"),
    @("05:15:38", "07:06", "215C_EL PATO", 111, "LP1912")
)

$r = 6
foreach ($row in $data2) {
    $ws2.Cells.Item($r, 1).Value = $row[0]
    $ws2.Cells.Item($r, 2).Value = $row[1]
    $ws2.Cells.Item($r, 3).Value = $row[2]
    $ws2.Cells.Item($r, 4).Value = $row[3]
    $ws2.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# Sheet "6203-6173": only the scrape timestamp changes (no data rows).
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 05:15:38"
